# impfmonitoring_three_notes.xlsx
#
# Commit: "Deal with unnamed notes columns (#11) - treat x columns as notes"
#
# On the "30.12.20" sheet (Worksheets item 2), three new unnamed/"note"
# columns (H and I) are introduced so that the extra per-state remarks
# ("blabla" / "another note" / "bar" / "test") can be captured alongside
# the existing A:G data table. This naturally extends the used range
# from A1:G23 to A1:I23 and adds four new shared strings.
#
# The sheet's view is also nudged: Excel re-derives topLeftCell/selection
# itself once new cells outside the previous frozen view are touched, so
# we just select the new last-entered note cell (H14) to mirror the
# author's final selection/cursor position recorded in the diff
# (activeCell="H14", no explicit topLeftCell override any more).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)   # "30.12.20"

# New "notes" cells (plain, unstyled text -- default style, like the
# other newly introduced note cells in this workbook).
$ws.Range("H2").Value  = "blabla"
$ws.Range("I2").Value  = "another note"
$ws.Range("I3").Value  = "bar"
$ws.Range("H14").Value = "test"

# Make the sheet active and park the selection on the last note cell
# that was added, matching the author's recorded cursor position.
[void]$ws.Activate()
[void]$ws.Range("H14").Select()
